$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Astro Resources")

$ws.Range("B11").Value = "SIMULATION LIST"
$ws.Range("C11").Value = "List of HTML, Native, and Java Astro SIMS"

$ws.Hyperlinks.Add(
    $ws.Cells.Item(11, 5),
    "https://docs.google.com/spreadsheets/d/1UN2LIh8TIoAYmp20fA9wYNC4XeJsuFte7fp4hPGfrMI/edit",
    "gid=0"
)

$ws.Range("E11").Value = "https://docs.google.com/spreadsheets/d/1UN2LIh8TIoAYmp20fA9wYNC4XeJsuFte7fp4hPGfrMI/edit#gid=0"
$ws.Range("E11").Style = "Hyperlink"

[void]$ws.Range("E12").Select()
